$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2721777.2
$ws.Range("J17").Value = 2809373.2
$ws.Range("L17").Value = 8428119.600000001
$ws.Range("N17").Value = -8428455.600000001

$ws.Range("H33").Value = 2925619
$ws.Range("I33").Value = 5848546
$ws.Range("K33").Value = 5848546
$ws.Range("M33").Value = -5848317

$ws.Range("H61").Value = 375.14285
$ws.Range("I61").Value = 375.14285
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1125.42855
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -953.4285500000001
$ws.Range("N61").ClearContents()

$ws.Range("H64").Value = 4111
$ws.Range("I64").Value = 4065.6667
$ws.Range("J64").Value = 4122.3335
$ws.Range("K64").Value = 4065.6667
$ws.Range("L64").Value = 4122.3335
$ws.Range("M64").Value = -3817.6667
$ws.Range("N64").Value = -4618.3335

$ws.Range("H67").Value = 4111
$ws.Range("I67").Value = 4065.6667
$ws.Range("J67").Value = 4122.3335
$ws.Range("K67").Value = 4065.6667
$ws.Range("L67").Value = 4122.3335
$ws.Range("M67").Value = -3207.6667
$ws.Range("N67").Value = -5838.3335

$ws.Range("H80").Value = 589688.25
$ws.Range("I80").Value = 1174.6666
$ws.Range("J80").Value = 910695.6
$ws.Range("K80").Value = 3523.9998
$ws.Range("L80").Value = 2732086.8
$ws.Range("M80").Value = -2525.9998
$ws.Range("N80").Value = -2734082.8

$ws.Range("H83").Value = 589688.25
$ws.Range("I83").Value = 1174.6666
$ws.Range("J83").Value = 910695.6
$ws.Range("K83").Value = 10571.9994
$ws.Range("L83").Value = 8196260.399999999
$ws.Range("M83").Value = -5579.999400000001
$ws.Range("N83").Value = -8206244.399999999

$ws.Range("H100").Value = 75205.71000000001
$ws.Range("I100").Value = 80606.16
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 80606.16
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -80065.16
$ws.Range("N100").Value = -6082

$ws.Range("H113").Value = 4529.0586
$ws.Range("I113").Value = 3250.4167
$ws.Range("J113").Value = 7597.8
$ws.Range("K113").Value = 3250.4167
$ws.Range("L113").Value = 7597.8
$ws.Range("M113").Value = 3.583299999999781
$ws.Range("N113").Value = -14105.8

$ws.Range("H132").Value = 736.4524
$ws.Range("I132").Value = 654.6389
$ws.Range("K132").Value = 1963.9167
$ws.Range("M132").Value = 566.0832999999998

$ws.Range("H135").Value = 4189.8335
$ws.Range("I135").Value = 4189.8335
$ws.Range("K135").Value = 37708.5015
$ws.Range("M135").Value = -35173.5015

$ws.Range("H137").Value = 9964.837
$ws.Range("I137").Value = 4335.174
$ws.Range("J137").Value = 14944.923
$ws.Range("K137").Value = 13005.522
$ws.Range("L137").Value = 44834.769
$ws.Range("M137").Value = -10455.522
$ws.Range("N137").Value = -49934.769

$ws.Range("H138").Value = 3440.5293
$ws.Range("I138").Value = 3331.8333
$ws.Range("J138").Value = 3499.818
$ws.Range("K138").Value = 9995.499899999999
$ws.Range("L138").Value = 10499.454
$ws.Range("M138").Value = -4855.499899999999
$ws.Range("N138").Value = -20779.454

$ws.Range("H141").Value = 3259.6
$ws.Range("I141").Value = 3666
$ws.Range("J141").Value = 2650
$ws.Range("K141").Value = 10998
$ws.Range("L141").Value = 7950
$ws.Range("M141").Value = -5818
$ws.Range("N141").Value = -18310

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 182.88889
$ws.Range("I4").Value = 280
$ws.Range("J4").Value = 170.75
$ws.Range("K4").Value = 280
$ws.Range("L4").Value = 170.75
$ws.Range("M4").Value = -164
$ws.Range("N4").Value = -402.75

$ws.Range("H32").Value = 2729.7234
$ws.Range("I32").Value = 2729.7234
$ws.Range("K32").Value = 2729.7234
$ws.Range("M32").Value = -2442.7234

$ws.Range("H61").Value = 9315.228999999999
$ws.Range("I61").Value = 7801.1
$ws.Range("K61").Value = 7801.1
$ws.Range("M61").Value = -7589.1

$ws.Range("H74").Value = 16857.334
$ws.Range("I74").Value = 17005.92
$ws.Range("K74").Value = 17005.92
$ws.Range("M74").Value = -16131.92

$ws.Range("H77").Value = 16857.334
$ws.Range("I77").Value = 17005.92
$ws.Range("K77").Value = 85029.59999999999
$ws.Range("M77").Value = -80661.59999999999

$ws.Range("H94").Value = 39999.5
$ws.Range("J94").Value = 39999.5
$ws.Range("L94").Value = 39999.5
$ws.Range("N94").Value = -41801.5

$ws.Range("H97").Value = 1358.3
$ws.Range("I97").Value = 1503.5
$ws.Range("K97").Value = 1503.5
$ws.Range("M97").Value = -1007.5

$ws.Range("H102").Value = 1322
$ws.Range("I102").Value = 1322
$ws.Range("K102").Value = 1322
$ws.Range("M102").Value = 300

$ws.Range("H114").Value = 82000
$ws.Range("J114").Value = 82000
$ws.Range("L114").Value = 82000
$ws.Range("N114").Value = -90678

$ws.Range("H122").Value = 2055.8667
$ws.Range("I122").Value = 2056.5833
$ws.Range("J122").Value = 2053
$ws.Range("K122").Value = 6169.749899999999
$ws.Range("L122").Value = 6159
$ws.Range("M122").Value = -3719.749899999999
$ws.Range("N122").Value = -11059

$ws.Range("H132").Value = 2420.743
$ws.Range("I132").Value = 1793.1666
$ws.Range("K132").Value = 5379.4998
$ws.Range("M132").Value = -2849.4998

$ws.Range("H136").Value = 9315.228999999999
$ws.Range("I136").Value = 7801.1
$ws.Range("K136").Value = 23403.3
$ws.Range("M136").Value = -20853.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 368.2857
$ws.Range("I64").Value = 166
$ws.Range("J64").Value = 402
$ws.Range("K64").Value = 166
$ws.Range("L64").Value = 402
$ws.Range("M64").Value = 59
$ws.Range("N64").Value = -852

$ws.Range("H67").Value = 368.2857
$ws.Range("I67").Value = 166
$ws.Range("J67").Value = 402
$ws.Range("K67").Value = 166
$ws.Range("L67").Value = 402
$ws.Range("M67").Value = 614
$ws.Range("N67").Value = -1962

$ws.Range("H94").Value = 3397.8572
$ws.Range("I94").Value = 2657.2
$ws.Range("J94").Value = 5249.5
$ws.Range("K94").Value = 2657.2
$ws.Range("L94").Value = 5249.5
$ws.Range("M94").Value = -2206.2
$ws.Range("N94").Value = -6151.5

$ws.Range("H105").Value = 4663.524
$ws.Range("I105").Value = 3964.8125
$ws.Range("K105").Value = 3964.8125
$ws.Range("M105").Value = -2217.8125

$ws.Range("H107").Value = 1358.6666
$ws.Range("I107").Value = 380.6
$ws.Range("J107").Value = 6249
$ws.Range("K107").Value = 380.6
$ws.Range("L107").Value = 6249
$ws.Range("M107").Value = 1539.4
$ws.Range("N107").Value = -10089

$ws.Range("H134").Value = 7372.741
$ws.Range("I134").Value = 3244.3076
$ws.Range("J134").Value = 15846.895
$ws.Range("K134").Value = 9732.9228
$ws.Range("L134").Value = 47540.685
$ws.Range("M134").Value = -7197.9228
$ws.Range("N134").Value = -52610.685

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5792.75
$ws.Range("I31").Value = 5786
$ws.Range("K31").Value = 5786
$ws.Range("M31").Value = -5491

$ws.Range("H34").Value = 5792.75
$ws.Range("I34").Value = 5786
$ws.Range("K34").Value = 5786
$ws.Range("M34").Value = -5584

$ws.Range("H58").Value = 5094.4
$ws.Range("I58").Value = 3492.5454
$ws.Range("K58").Value = 3492.5454
$ws.Range("M58").Value = -3289.5454

$ws.Range("H62").Value = 206896.4
$ws.Range("I62").Value = 502250
$ws.Range("J62").Value = 9994
$ws.Range("K62").Value = 502250
$ws.Range("L62").Value = 9994
$ws.Range("M62").Value = -501626
$ws.Range("N62").Value = -11242

$ws.Range("H65").Value = 206896.4
$ws.Range("I65").Value = 502250
$ws.Range("J65").Value = 9994
$ws.Range("K65").Value = 2511250
$ws.Range("L65").Value = 49970
$ws.Range("M65").Value = -2508130
$ws.Range("N65").Value = -56210

$ws.Range("H86").Value = 3834.5293
$ws.Range("I86").Value = 2876
$ws.Range("J86").Value = 4505.5
$ws.Range("K86").Value = 2876
$ws.Range("L86").Value = 4505.5
$ws.Range("M86").Value = -1753
$ws.Range("N86").Value = -6751.5

$ws.Range("H89").Value = 3834.5293
$ws.Range("I89").Value = 2876
$ws.Range("J89").Value = 4505.5
$ws.Range("K89").Value = 14380
$ws.Range("L89").Value = 22527.5
$ws.Range("M89").Value = -8764
$ws.Range("N89").Value = -33759.5

$ws.Range("H99").Value = 11690.167
$ws.Range("I99").Value = 6399.2
$ws.Range("K99").Value = 6399.2
$ws.Range("M99").Value = -4901.2

$ws.Range("H122").Value = 3272.1
$ws.Range("I122").Value = 2523
$ws.Range("K122").Value = 7569
$ws.Range("M122").Value = -5119

$ws.Range("H126").Value = 11690.167
$ws.Range("I126").Value = 6399.2
$ws.Range("K126").Value = 19197.6
$ws.Range("M126").Value = -16727.6

$ws.Range("H132").Value = 30626.88
$ws.Range("I132").Value = 20339.217
$ws.Range("K132").Value = 61017.651
$ws.Range("M132").Value = -58487.651

$ws.Range("H134").Value = 5999.933
$ws.Range("I134").Value = 3888.2
$ws.Range("J134").Value = 10223.4
$ws.Range("K134").Value = 11664.6
$ws.Range("L134").Value = 30670.2
$ws.Range("M134").Value = -9129.599999999999
$ws.Range("N134").Value = -35740.2

$ws.Range("H136").Value = 5094.4
$ws.Range("I136").Value = 3492.5454
$ws.Range("K136").Value = 10477.6362
$ws.Range("M136").Value = -7927.636200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1868.5
$ws.Range("I5").Value = 1197.8334
$ws.Range("J5").Value = 2092.0557
$ws.Range("K5").Value = 3593.5002
$ws.Range("L5").Value = 6276.1671
$ws.Range("M5").Value = -3481.5002
$ws.Range("N5").Value = -6500.1671

$ws.Range("H11").Value = 956.5
$ws.Range("J11").Value = 529.3333
$ws.Range("L11").Value = 1587.9999
$ws.Range("N11").Value = -1867.9999

$ws.Range("H38").Value = 2466.0715
$ws.Range("I38").Value = 485.8889
$ws.Range("J38").Value = 6030.4
$ws.Range("K38").Value = 1457.6667
$ws.Range("L38").Value = 18091.2
$ws.Range("M38").Value = -1110.6667
$ws.Range("N38").Value = -18785.2

$ws.Range("H59").Value = 1300
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 2068.1667
$ws.Range("I60").Value = 1990
$ws.Range("J60").Value = 2224.5
$ws.Range("K60").Value = 5970
$ws.Range("L60").Value = 6673.5
$ws.Range("M60").Value = -5719
$ws.Range("N60").Value = -7175.5

$ws.Range("H61").Value = 69.59999999999999
$ws.Range("I61").Value = 69.59999999999999
$ws.Range("K61").Value = 208.8
$ws.Range("M61").Value = 6.200000000000017

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H63").Value = 2998.5
$ws.Range("I63").Value = 2998.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 8995.5
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -8246.5
$ws.Range("N63").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H66").Value = 2998.5
$ws.Range("I66").Value = 2998.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 26986.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -23242.5
$ws.Range("N66").ClearContents()

$ws.Range("H69").Value = 11749.667
$ws.Range("J69").Value = 12499.8
$ws.Range("L69").Value = 37499.39999999999
$ws.Range("N69").Value = -39121.39999999999

$ws.Range("H70").Value = 1000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H72").Value = 11749.667
$ws.Range("J72").Value = 12499.8
$ws.Range("L72").Value = 112498.2
$ws.Range("N72").Value = -120610.2

$ws.Range("H73").Value = 1000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H74").Value = 16000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 16000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 48000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -50122

$ws.Range("H75").Value = 111111784
$ws.Range("J75").Value = 111111784
$ws.Range("L75").Value = 333335352
$ws.Range("N75").Value = -333337348

$ws.Range("H77").Value = 16000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 16000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 144000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -154608

$ws.Range("H78").Value = 111111784
$ws.Range("J78").Value = 111111784
$ws.Range("L78").Value = 1000006056
$ws.Range("N78").Value = -1000016040

$ws.Range("H80").Value = 54999.5
$ws.Range("I80").Value = 10000
$ws.Range("J80").Value = 99999
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 299997
$ws.Range("M80").Value = -29064
$ws.Range("N80").Value = -301869

$ws.Range("H83").Value = 54999.5
$ws.Range("I83").Value = 10000
$ws.Range("J83").Value = 99999
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 899991
$ws.Range("M83").Value = -85320
$ws.Range("N83").Value = -909351

$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws.Range("H97").Value = 31130
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 31130
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 93390
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -94382

$ws.Range("H121").Value = 2173664
$ws.Range("I121").Value = 2248.7273
$ws.Range("K121").Value = 6746.1819
$ws.Range("M121").Value = -5436.1819

$ws.Range("H135").Value = 1868.5
$ws.Range("I135").Value = 1197.8334
$ws.Range("J135").Value = 2092.0557
$ws.Range("K135").Value = 10780.5006
$ws.Range("L135").Value = 18828.5013
$ws.Range("M135").Value = -8245.500599999999
$ws.Range("N135").Value = -23898.5013

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 19999.334
$ws.Range("I5").Value = 15000
$ws.Range("J5").Value = 24998.666
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 24998.666
$ws.Range("M5").Value = -14888
$ws.Range("N5").Value = -25222.666

$ws.Range("H97").Value = 1290.5
$ws.Range("I97").Value = 1417.6666
$ws.Range("J97").Value = 1163.3334
$ws.Range("K97").Value = 1417.6666
$ws.Range("L97").Value = 1163.3334
$ws.Range("M97").Value = -921.6666
$ws.Range("N97").Value = -2155.3334

$ws.Range("H102").Value = 1159.5238
$ws.Range("I102").Value = 1185.5294
$ws.Range("K102").Value = 1185.5294
$ws.Range("M102").Value = 436.4706000000001

$ws.Range("H123").Value = 41666.668
$ws.Range("J123").Value = 41666.668
$ws.Range("L123").Value = 41666.668
$ws.Range("N123").Value = -46566.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10999.8
$ws.Range("I7").Value = 9999
$ws.Range("J7").Value = 11250
$ws.Range("K7").Value = 9999
$ws.Range("L7").Value = 11250
$ws.Range("M7").Value = -9887
$ws.Range("N7").Value = -11474

$ws.Range("H22").Value = 2283.7778
$ws.Range("J22").Value = 2463.6924
$ws.Range("L22").Value = 2463.6924
$ws.Range("N22").Value = -3053.6924

$ws.Range("H27").Value = 2283.7778
$ws.Range("J27").Value = 2463.6924
$ws.Range("L27").Value = 2463.6924
$ws.Range("N27").Value = -2677.6924

$ws.Range("H61").Value = 3821.6667
$ws.Range("I61").Value = 2399.5
$ws.Range("K61").Value = 2399.5
$ws.Range("M61").Value = -2197.5

$ws.Range("H100").Value = 5028.1665
$ws.Range("I100").Value = 3166.3333
$ws.Range("K100").Value = 3166.3333
$ws.Range("M100").Value = -2625.3333

$ws.Range("H113").Value = 3821.6667
$ws.Range("I113").Value = 2399.5
$ws.Range("K113").Value = 2399.5
$ws.Range("M113").Value = -229.5

$ws.Range("H122").Value = 4699.143
$ws.Range("J122").Value = 3565.6667
$ws.Range("L122").Value = 10697.0001
$ws.Range("N122").Value = -15597.0001

$ws.Range("H126").Value = 10999.8
$ws.Range("I126").Value = 9999
$ws.Range("J126").Value = 11250
$ws.Range("K126").Value = 29997
$ws.Range("L126").Value = 33750
$ws.Range("M126").Value = -27527
$ws.Range("N126").Value = -38690

$ws.Range("H132").Value = 4629.727
$ws.Range("I132").Value = 4646.303
$ws.Range("J132").Value = 4580
$ws.Range("K132").Value = 13938.909
$ws.Range("L132").Value = 13740
$ws.Range("M132").Value = -11408.909
$ws.Range("N132").Value = -18800

$ws.Range("H136").Value = 3127.7036
$ws.Range("I136").Value = 2597.8408
$ws.Range("J136").Value = 5459.1
$ws.Range("K136").Value = 7793.5224
$ws.Range("L136").Value = 16377.3
$ws.Range("M136").Value = -5243.5224
$ws.Range("N136").Value = -21477.3

$ws.Range("H137").Value = 59000
$ws.Range("J137").Value = 59000
$ws.Range("L137").Value = 59000
$ws.Range("N137").Value = -69200

$ws.Range("H140").Value = 76652.86
$ws.Range("J140").Value = 76652.86
$ws.Range("L140").Value = 76652.86
$ws.Range("N140").Value = -87012.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2999.5
$ws.Range("I2").Value = 2999.5
$ws.Range("K2").Value = 2999.5
$ws.Range("M2").Value = -2887.5

$ws.Range("H74").Value = 4866.5713
$ws.Range("J74").Value = 3380.2
$ws.Range("L74").Value = 3380.2
$ws.Range("N74").Value = -5252.2

$ws.Range("H77").Value = 4866.5713
$ws.Range("J77").Value = 3380.2
$ws.Range("L77").Value = 10140.6
$ws.Range("N77").Value = -19500.6

$ws.Range("H122").Value = 7993.5713
$ws.Range("I122").Value = 5265.636
$ws.Range("K122").Value = 15796.908
$ws.Range("M122").Value = -13346.908

$ws.Range("H132").Value = 164813.48
$ws.Range("I132").Value = 267657.78
$ws.Range("K132").Value = 802973.3400000001
$ws.Range("M132").Value = -800443.3400000001

$ws.Range("H136").Value = 5002555
$ws.Range("I136").Value = 10001954
$ws.Range("K136").Value = 30005862
$ws.Range("M136").Value = -30003312
